$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.504
$ws.Range("D7").Value = -7.445
$ws.Range("C10").Value = -12.937
$ws.Range("C12").Value = -11.066
$ws.Range("D15").Value = -8.378000000000002
$ws.Range("C18").Value = -13.265
$ws.Range("E18").Value = 16.974
$ws.Range("E19").Value = 16.606
$ws.Range("D20").Value = -7.517
$ws.Range("E27").Value = 16.385
$ws.Range("D29").Value = -7.292
$ws.Range("D30").Value = -7.219999999999999
$ws.Range("D31").Value = -8.096
$ws.Range("C37").Value = -13.475
$ws.Range("D40").Value = -7.641
$ws.Range("E42").Value = 16.555
$ws.Range("E44").Value = 16.679
$ws.Range("E47").Value = 16.467
$ws.Range("C55").Value = -13.916
$ws.Range("E58").Value = 16.596
$ws.Range("C68").Value = -11.167
$ws.Range("D68").Value = -6.879
$ws.Range("E73").Value = 16.667
$ws.Range("D76").Value = -7.311999999999999
$ws.Range("C77").Value = -13.117
$ws.Range("C78").Value = -13.214
$ws.Range("D87").Value = -8.371
$ws.Range("D88").Value = -7.973999999999999
$ws.Range("E95").Value = 17.397
$ws.Range("D96").Value = -7.267
$ws.Range("D98").Value = -8.228
$ws.Range("D101").Value = -7.616
$ws.Range("E101").Value = 16.44
$ws.Range("D102").Value = -8.036
